
# The diff touches xl/worksheets/sheet1.xml, which the workbook relationships
# map to the sheet named "R1". Grab it explicitly (rather than relying on
# ActiveSheet) so the edit lands on the right tab regardless of selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("R1")

# Append a new outage row (row 7) below the existing data (rows 1-6),
# matching the columns: PCM, Region, Count sites, Hub Site, Fault Level,
# PCM Created At, Elapsed Duration(Hrs), Creat Fault First Time Occured,
# Power Source, Battery Backup Status, EM Field Feedback, Site Owner.
$row = 7
$ws.Cells.Item($row, 1).Value = ""          # A7 - PCM
$ws.Cells.Item($row, 2).Value = "R4"        # B7 - Region
$ws.Cells.Item($row, 3).Value = ""          # C7 - Count sites
$ws.Cells.Item($row, 4).Value = "LTH2121"   # D7 - Hub Site
$ws.Cells.Item($row, 5).Value = ""          # E7 - Fault Level
$ws.Cells.Item($row, 6).Value = ""          # F7 - PCM Created At
$ws.Cells.Item($row, 7).Value = ""          # G7 - Elapsed Duration(Hrs)
$ws.Cells.Item($row, 8).Value = ""          # H7 - Creat Fault First Time Occured
$ws.Cells.Item($row, 9).Value = "SCECO"     # I7 - Power Source
$ws.Cells.Item($row, 10).Value = "Dead"     # J7 - Battery Backup Status
$ws.Cells.Item($row, 11).Value = ""         # K7 - EM Field Feedback
$ws.Cells.Item($row, 12).Value = "Latis"    # L7 - Site Owner
